# Fruta / hortaliza, semanal
# Insert a new daily observation row at row 44 (pushing the existing
# rows 44-124 down to 45-125) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(44).Insert()

$ws.Range("A44").Value = 6
$ws.Range("B44").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C44").Value = "Metropolitana"
$ws.Range("D44").Value = 44540
$ws.Range("E44").Value = 13
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100101
$ws.Range("H44").Value = "Berries"
$ws.Range("I44").Value = 100101004
$ws.Range("J44").Value = "Frambuesa"
$ws.Range("K44").Value = "Sin especificar"
$ws.Range("L44").Value = "Especial"
$ws.Range("M44").Value = 500
$ws.Range("N44").Value = 8000
$ws.Range("O44").Value = 8000
$ws.Range("P44").Value = 8000
$ws.Range("Q44").Value = "$/bandeja 2 kilos"
$ws.Range("R44").Value = "Provincia de Curicó"
$ws.Range("S44").Value = 4000
$ws.Range("T44").Value = 2
